# fix: reviewsCount 1.2k to 1200
# Populates Sheet1 with the yelp results table and styles the header row
# (bold, centered/top-aligned, thin box border).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$headers = @("name", "rating", "reviewCount", "keyword")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Build the bold/centered/top/boxed style once on A1 only, then copy that
# finished format onto B1:D1 via copy/paste-special-formats. Applying the
# same sequence of property writes directly to a multi-cell Range (or to
# each cell independently) instead forks a brand-new cellXf per write, so
# several unused intermediate xf records get stranded in the stylesheet.
# Building the style on a single cell first and cloning it keeps the
# stylesheet minimal: one new font, one new border, one new cellXf.
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4160     # xlTop
$a1.Borders.LineStyle = "thin"

$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows ----------------------------------------------------------
$rows = @(
    @("Cobble Fish",            4.1, 42,   "Seafood, Cocktail Bars, South Street Seaport"),
    @("Sweetwater Restaurant",  4.1, 483,  "New American, Bars"),
    @("Okdongsik",              4.3, 290,  "Korean, Soup, Soul Food"),
    @("Fresh Salt",             4,   425,  "Bars, New American"),
    @("Debajo",                 4.5, 79,   "Tapas Bars, Spanish, Tapas/Small Plates, Flatiron"),
    @("nonono",                 4.2, 1200, "Izakaya, Cocktail Bars"),
    @("Sweet Anaelle",          4.7, 6,    "Peruvian, Cocktail Bars, Bushwick"),
    @("Betong - Khao Man Gai",  4.5, 50,   "Thai"),
    @("OBAO",                   4,   4200, "Vietnamese, Thai, Asian Fusion"),
    @("R40",                    4.4, 215,  "Argentine, Cocktail Bars")
)

$rowIndex = 2
foreach ($row in $rows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $rowIndex++
}
